$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 ("Videos" / "Put individual links to videos on YouTube on Video pages")
$ws.Rows.Item(3).Delete()

# After the above deletion, the two rows that used to be 13 and 14
# ("Cleanup" / "Remove code that's not used any more" and
#  "" / "Remove content that is no longer used") are now rows 12 and 13.
$ws.Range("A12:B13").EntireRow.Delete()

# Update the active selection to match the saved view state
$ws.Range("A16").Select()
